# Add KPI_History tab with columns
# - Rename the default sheet to "KPI_History"
# - Write the 14 KPI column headers into row 1 (this also seeds the
#   shared-strings table)
# - Bold the header row, with a text ("@") number format on the first
#   column (quarter), matching the authored workbook
# - Turn the header row into a real Excel Table ("Table_1") spanning
#   A1:N1 with banded rows / highlighted first+last column and no
#   auto-filter dropdowns, using a custom table style name

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename Sheet1 -> KPI_History
$ws.Name = "KPI_History"

# 2. Column headers (order matters - maps to columns A..N)
$headers = @(
    "quarter",
    "revenue",
    "yoy_growth_%",
    "gross_margin_%",
    "s&m_%rev",
    "r&d_%rev",
    "g&a_%rev",
    "arr",
    "nrr_%",
    "churn_%",
    "arpu",
    "cfo",
    "capex",
    "fcf"
)

for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 3. Header formatting: bold everywhere, text format on the quarter column
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").NumberFormat = "@"
$ws.Range("B1:N1").Font.Bold = $true

# 4. Turn A1:N1 into a Table
$headerRange = $ws.Range("A1:N1")
$tbl = $ws.ListObjects.Add(1, $headerRange, 0, 1)
$tbl.Name = "Table_1"

$tbl.ShowAutoFilter = $false
$tbl.ShowTableStyleFirstColumn = $true
$tbl.ShowTableStyleLastColumn = $true
$tbl.ShowTableStyleRowStripes = $true
$tbl.ShowTableStyleColumnStripes = $false
$tbl.TableStyle = "KPI_History-style"
